$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Ручки КПП" (3rd sheet): rename "spaco"/"SPACO" -> "sparco"/"SPARCO"
# in the gear-shift-knob item rows, and bump row 3's height.
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# (new shared-string entries are created in the order the cells are
# written, so write all of column A first, then all of column B, to
# reproduce the original authoring order)
$ws3.Range("A2").Value2 = "Ручка переключения передач sparco (сереб-черн-кр)"
$ws3.Range("A3").Value2 = "Ручка переключения передач sparco (сереб - черн)"
$ws3.Range("A4").Value2 = "Ручка переключения передач sparco (черн - бел)"
$ws3.Range("A5").Value2 = "Ручка переключения передач sparco (черн - бел - кр)"
$ws3.Range("A6").Value2 = "Ручка переключения передач sparco (черн-кр)"

$ws3.Range("B2").Value2 = "Ручка переключения кпп SPARCO (сереб - черн - кр), в комплекте переходные гайки под разную резьбу штока - M8/М10/M12"
$ws3.Range("B3").Value2 = "Ручка переключения кпп SPARCO (сереб - черн), в комплекте переходные гайки под разную резьбу штока - M8/М10/M12"
$ws3.Range("B4").Value2 = "Ручка переключения кпп SPARCO (черн - бел), в комплекте переходные гайки под разную резьбу штока - M8/М10/M12"
$ws3.Range("B5").Value2 = "Ручка переключения кпп SPARCO (черн - бел - кр), в комплекте переходные гайки под разную резьбу штока - M8/М10/M12"
$ws3.Range("B6").Value2 = "Ручка переключения кпп SPARCO (черн-кр), в комплекте переходные гайки под разную резьбу штока - M8/М10/M12"

# Row 3 grew taller once its label text got longer.
$ws3.Rows.Item(3).RowHeight = 45

# Cursor/selection left on G17 when the edit was saved.
$ws3.Range("G17").Select()

# -----------------------------------------------------------------
# Sheet "Накладки на педали" (1st sheet): re-enter the percent/diff
# formulas over the whole column so they become shared formulas,
# and correct the cost for the last item which lowers its margin.
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2:F11").Formula = "=(E2-D2)/D2*100"
$ws1.Range("G2:G11").Formula = "=E2-D2"

$ws1.Range("D11").Value2 = 1021.96

# Cursor/selection left on D12 when the edit was saved.
$ws1.Range("D12").Select()
